$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0) The trailing "_GoBack" bookmark currently sits at the end of the
#    "Celular: ..." paragraph. The edit relocates it onto the new
#    "Enlace para ingresar al github:" paragraph, so drop the old one
#    here and we'll re-create it in the right spot below.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 1) Mint the two external hyperlink relationships through the real
#    Word object model (Hyperlinks.Add both registers the
#    relationship in word/_rels/document.xml.rels AND drops a
#    <w:hyperlink> into the body). We don't hard-code the rXX ids the
#    engine allocates - we read them back after the fact.
# ------------------------------------------------------------------
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$tail.Collapse(0)
$d.Hyperlinks.Add($tail, "https://github.com/JoelAlBe", $null, $null, "placeholder1") | Out-Null

$tail2 = $d.Content
$tail2.Collapse(0)
$tail2.InsertParagraphAfter()
$tail2.Collapse(0)
$d.Hyperlinks.Add($tail2, "https://trello.com/b/NSSjNTkK/gestion-de-historiales", $null, $null, "placeholder2") | Out-Null

# Recover the ids minted by the two Hyperlinks.Add calls above by
# reading the raw markup of the two placeholder paragraphs back out
# (WordOpenXML always serialises the whole package, so walk every
# <w:hyperlink r:id="..."> occurrence and keep the last two - those
# are necessarily the ones we just minted, in call order).
$lastTwoStart = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range.Start
$placeholderRange = $d.Range($lastTwoStart, $d.Content.End)
$rawXml = $placeholderRange.WordOpenXML

$foundIds = @()
$rest = $rawXml
while ($rest -match 'w:hyperlink r:id="(rId\d+)"') {
  $foundIds += $matches[1]
  $idx = $rest.IndexOf($matches[0]) + $matches[0].Length
  $rest = $rest.Substring($idx)
}
$n = $foundIds.Count
$githubRid = $foundIds[$n - 2]
$trelloRid = $foundIds[$n - 1]

# ------------------------------------------------------------------
# 2) Replace the two placeholder paragraphs with the exact OOXML the
#    commit produced: custom run-splitting for the spell-checker's
#    proofErr markers around "JoelAlBe" / "Gestion", and the bookmark
#    relocated onto the github paragraph. InsertXML *replaces* the
#    contents of the range it is called on.
# ------------------------------------------------------------------
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

$newXml = @"
<w:p $ns>
  <w:pPr>
    <w:rPr>
      <w:lang w:eastAsia="es-MX"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:shd w:val="clear" w:color="auto" w:fill="8EAADB" w:themeFill="accent1" w:themeFillTint="99"/>
      <w:lang w:eastAsia="es-MX"/>
    </w:rPr>
    <w:t>Enlace para ingresar al github:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:eastAsia="es-MX"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:hyperlink r:id="$githubRid" w:history="1">
    <w:proofErr w:type="spellStart"/>
    <w:r>
      <w:t>JoelAlBe</w:t>
    </w:r>
    <w:proofErr w:type="spellEnd"/>
    <w:r>
      <w:t xml:space="preserve"> (github.com)</w:t>
    </w:r>
  </w:hyperlink>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:rPr>
      <w:lang w:eastAsia="es-MX"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:shd w:val="clear" w:color="auto" w:fill="8EAADB" w:themeFill="accent1" w:themeFillTint="99"/>
      <w:lang w:eastAsia="es-MX"/>
    </w:rPr>
    <w:t>Enlace para ingresar a trello:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:eastAsia="es-MX"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:eastAsia="es-MX"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:hyperlink r:id="$trelloRid" w:history="1">
    <w:proofErr w:type="spellStart"/>
    <w:r>
      <w:t>Gestion</w:t>
    </w:r>
    <w:proofErr w:type="spellEnd"/>
    <w:r>
      <w:t xml:space="preserve"> de historiales | Trello</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
"@

$placeholderRange.InsertXML($newXml)

# ------------------------------------------------------------------
# 3) InsertXML drops character-style references (<w:rStyle>), so
#    reapply the "Hipervnculo" hyperlink character style to the
#    display text of both hyperlinks we just inserted, exactly as
#    Word itself would when the hyperlink was typed/created.
# ------------------------------------------------------------------
$hCount = $d.Hyperlinks.Count
$d.Hyperlinks.Item($hCount - 1).Range.Style = "Hipervnculo"
$d.Hyperlinks.Item($hCount).Range.Style = "Hipervnculo"

Write-Output "done"
